# Correction in Excel file:
# The "N'" column (E) values were off; correct them to the actual rounded
# N values used for the verification column (F), which recalculates
# automatically from the formulas already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 840
$ws.Range("E3").Value = 1300
$ws.Range("E4").Value = 1761

# Move the active selection to E3, matching where the author left off.
$ws.Range("E3").Select()
